$wb = $excel.ActiveWorkbook

# Rename the "ATT" sheet to "FilterOptions"
$ws = $wb.Worksheets.Item("ATT")
$ws.Name = "FilterOptions"

# Populate the new filter-options list (2 more test cases worth of reference data)
$values = @(
    "Optios",
    "Drive-thru ATM",
    "Walk-up ATM",
    "Cardless ATM",
    "ATM accepts deposits and credit card payments",
    "ATM cash withdrawals only (deposits not accepted)",
    "ATM located inside",
    "Dedicated Business Teller",
    "Glass barrier at Customer Service desk",
    "Drive-thru Teller Services",
    "Video Conferencing",
    "Express financial center",
    "Advanced Center™ with Video Chat",
    "Accepts appointments",
    "Financial Solutions Advisor",
    "Home Loans Specialist",
    "Centralized Small Business Banker",
    "Notary",
    "Commercial Deposits Accepted",
    "Night Deposit Accepted",
    "Change Orders",
    "ATM Services Available",
    "Open Saturdays",
    "Open Sundays",
    "Currently open"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Widen column A so the longer labels are readable
$ws.Columns.Item(1).ColumnWidth = 38.333333333333

# Make FilterOptions the active sheet/tab and select A7 (matches saved view state)
$ws.Activate()
$ws.Range("A7").Select()
